$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Items (column D) cell text: lower-case the item names ---
$ws.Range("D3").Value = "gun%100"
$ws.Range("D4").Value = "batteries%100"
$ws.Range("D8").Value = "knife%100"

# --- Update Nodes (column E) cell text ---
$ws.Range("E7").Value = "escapepod%100%false"

# --- Replace short room names in Description column with full descriptions ---
$ws.Range("C7").Value = "Pod room is where the escape pod is located with the help of which we can leave the spaceship."
$ws.Range("C5").Value = "Medbay is the room with all the necessary medical supplies."
$ws.Range("C6").Value = "Lab is where astronauts can perform various tests."
$ws.Range("C8").Value = "Cafeteria is where we all the food is present."
$ws.Range("C9").Value = "Gym has all the tools for staying fit while on the spaceship"
$ws.Range("C10").Value = "Cryogenics is where we perform research on life expectancy."

# --- Widen column E and move the active selection, matching the author's last UI state ---
$ws.Columns.Item(5).ColumnWidth = 21.5
$ws.Range("C16").Select()
